# Update row 9 (Ano 2025) values in the faturamento_anual sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 2831438.47
$ws.Range("C9").Value = 445261.99
$ws.Range("D9").Value = 3276700.46
$ws.Range("E9").Value = 13.5887303534605
$ws.Range("F9").Value = 86.41126964653949
$ws.Range("G9").Value = -56.96764884513605
$ws.Range("H9").Value = -48.86813495821934
$ws.Range("I9").Value = 28076
$ws.Range("J9").Value = 1199
$ws.Range("K9").Value = 29275
$ws.Range("L9").Value = 20177
$ws.Range("M9").Value = 162.3978024483323
$ws.Range("N9").Value = 10.87228740597346
